# Fruta / hortaliza, semanal
# Insert 3 new weekly data rows at the top of the data block (rows 221-223),
# pushing the existing rows 221-289 down to 224-292.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 221 (shifts existing rows 221:289 -> 224:292)
$ws.Range("A221:T223").EntireRow.Insert()

# New row 221
$ws.Range("A221").Value = 4
$ws.Range("B221").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C221").Value = "Los Lagos"
$ws.Range("D221").Value = 44876
$ws.Range("E221").Value = 10
$ws.Range("F221").Value = "Fruta"
$ws.Range("G221").Value = 100101
$ws.Range("H221").Value = "Berries"
$ws.Range("I221").Value = 100112025
$ws.Range("J221").Value = "Frutilla"
$ws.Range("K221").Value = "Sin especificar"
$ws.Range("L221").Value = "Especial"
$ws.Range("M221").Value = 200
$ws.Range("N221").Value = 13000
$ws.Range("O221").Value = 13000
$ws.Range("P221").Value = 13000
$ws.Range("Q221").Value = "$/bandeja 7 kilos"
$ws.Range("R221").Value = "Provincia de Melipilla"
$ws.Range("S221").Value = 1857
$ws.Range("T221").Value = 7

# New row 222
$ws.Range("A222").Value = 4
$ws.Range("B222").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C222").Value = "Los Lagos"
$ws.Range("D222").Value = 44876
$ws.Range("E222").Value = 10
$ws.Range("F222").Value = "Fruta"
$ws.Range("G222").Value = 100101
$ws.Range("H222").Value = "Berries"
$ws.Range("I222").Value = 100112025
$ws.Range("J222").Value = "Frutilla"
$ws.Range("K222").Value = "Sin especificar"
$ws.Range("L222").Value = "Primera"
$ws.Range("M222").Value = 500
$ws.Range("N222").Value = 10000
$ws.Range("O222").Value = 11000
$ws.Range("P222").Value = 10500
$ws.Range("Q222").Value = "$/bandeja 7 kilos"
$ws.Range("R222").Value = "Provincia de Melipilla"
$ws.Range("S222").Value = 1500
$ws.Range("T222").Value = 7

# New row 223
$ws.Range("A223").Value = 4
$ws.Range("B223").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C223").Value = "Los Lagos"
$ws.Range("D223").Value = 44876
$ws.Range("E223").Value = 10
$ws.Range("F223").Value = "Fruta"
$ws.Range("G223").Value = 100101
$ws.Range("H223").Value = "Berries"
$ws.Range("I223").Value = 100112025
$ws.Range("J223").Value = "Frutilla"
$ws.Range("K223").Value = "Sin especificar"
$ws.Range("L223").Value = "Primera"
$ws.Range("M223").Value = 600
$ws.Range("N223").Value = 10000
$ws.Range("O223").Value = 11000
$ws.Range("P223").Value = 10500
$ws.Range("Q223").Value = "$/caja 7 kilos"
$ws.Range("R223").Value = "Región de La Araucanía"
$ws.Range("S223").Value = 1500
$ws.Range("T223").Value = 7
